# Insert a new data row at row 434 (pushing the existing rows 434..535 down
# to 435..536) and populate it with the new record, matching the diff:
#   dimension A1:R535 -> A1:R536
#   new row 434: 2023-10-12, Rodeo, "1a (cosecha lavada)", 150, 25000, 25000,
#                25000, "$/malla 25 kilos", "Región de La Araucanía", 1000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("434").Insert()

$ws.Range("A434").Value = 11
$ws.Range("B434").Value = "Vega Monumental Concepción"
$ws.Range("C434").Value = "Bíobío"
$ws.Range("D434").Value = 45211
$ws.Range("E434").Value = 8
$ws.Range("F434").Value = 100114001
$ws.Range("G434").Value = "Papa"
$ws.Range("H434").Value = "Rodeo"
$ws.Range("I434").Value = "1a (cosecha lavada)"
$ws.Range("J434").Value = 150
$ws.Range("K434").Value = 25000
$ws.Range("L434").Value = 25000
$ws.Range("M434").Value = 25000
$ws.Range("N434").Value = "`$/malla 25 kilos"
$ws.Range("O434").Value = "Región de La Araucanía"
$ws.Range("P434").Value = 1000
$ws.Range("Q434").Value = 25
$ws.Range("R434").Value = "Hortaliza"
